# update docker & change logic crawl data
# Replace the crawled menu items with the freshly-scraped set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Row 2: Má Đùi gà chiên sốt tiêu đen
$ws.Range("A2").Value = "Má Đùi gà chiên sốt tiêu đen"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr5qnme94k9l6e"

# Row 3: Thịt kho đậu hũ
$ws.Range("A3").Value = "Thịt kho đậu hũ"
$ws.Range("B3").Value = "• Cơm gạo dẻo, nhiều rau xanh" + $nl + "• Có canh ăn kèm miễn phí" + $nl + "• Hộp đựng 4 ngăn tiện lợi"
$ws.Range("C3").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lqxef49eu3t09d"

# Row 4: Đậu hũ chiên sốt thịt bằm
$ws.Range("A4").Value = "Đậu hũ chiên sốt thịt bằm"
$ws.Range("B4").Value = "• Cơm gạo dẻo, nhiều rau xanh" + $nl + "• Có canh ăn kemf"
$ws.Range("C4").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr5fjziuklnd99"

# Row 5: Cá nục kho cải chua
$ws.Range("A5").Value = "Cá nục kho cải chua"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lqy2h5xy2oxw51"

# Row 6: Vịt kho gừng
$ws.Range("A6").Value = "Vịt kho gừng"
$ws.Range("B6").Value = "• Cơm gạo dẻo, nhiều rau xanh"
$ws.Range("C6").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lqxo3pgolk1wd7"

# Row 7: Canh chả cá thác lác khổ hoa bào
$ws.Range("A7").Value = "Canh chả cá thác lác khổ hoa bào"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr30e7y79nyxb7"

# Row 8: Bánh canh sườn non, trứng cút
$ws.Range("A8").Value = "Bánh canh sườn non, trứng cút"
$ws.Range("B8").Value = "• Cá lóc phi lê xé nhỏ ăn kèm chả ram chiên giòn" + $nl + "• Giá đã bao gồm đủ đồ xào và canh ăn kèm." + $nl + "• Rất vui được ăn trưa cùng bạn."
$ws.Range("C8").Value = "https://mms.img.susercontent.com/vn-11134517-7r98o-lr5zplxvkslgf0"

# Old row 9 (Canh khoai sọ hầm xương) is gone entirely; remaining
# "extra item" rows (Cơm thêm, Trứng ốp la, Cải chua trộn, Đồ ăn thêm,
# Đồ ăn thêm bò lúc lắc) shift up by one row to close the gap.
$ws.Rows.Item(9).Delete()
